# Add a new attendance column (BP) for the latest training date (2025-10-22,
# serial 45952), mirroring the existing BO column's layout.
#
# We set all the new cell *values* first (so the workbook's dependency graph
# correctly marks the summary formulas in columns B/C/F/I as dirty and they
# recompute against the newly-widened used range), and only afterwards copy
# the visual formatting (number format + alignment) from column BO onto the
# new column BP. Doing the formatting copy before the value writes causes
# the summary formulas to keep their stale cached results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: new training-session date for column BP, row 1.
$ws.Range("BP1").Value = 45952

# Per-player attendance mark for the new session (column BP), rows 2-29.
# Row 12's player had already left the team before this column range, so it
# is intentionally skipped (it has no entry in column BO either).
$ws.Range("BP2").Value = "P"
$ws.Range("BP3").Value = "P"
$ws.Range("BP4").Value = "P"
$ws.Range("BP5").Value = "P"
$ws.Range("BP6").Value = "B"
$ws.Range("BP7").Value = "P"
$ws.Range("BP8").Value = "B"
$ws.Range("BP9").Value = "P"
$ws.Range("BP10").Value = "P"
$ws.Range("BP11").Value = "P"
$ws.Range("BP13").Value = "B"
$ws.Range("BP14").Value = "P"
$ws.Range("BP15").Value = "P"
$ws.Range("BP16").Value = "B"
$ws.Range("BP17").Value = "RH"
$ws.Range("BP18").Value = "P"
$ws.Range("BP19").Value = "B"
$ws.Range("BP20").Value = "P"
$ws.Range("BP21").Value = "B"
$ws.Range("BP22").Value = "P"
$ws.Range("BP23").Value = "B"
$ws.Range("BP24").Value = "B"
$ws.Range("BP25").Value = "P"
$ws.Range("BP26").Value = "P"
$ws.Range("BP27").Value = "P"
$ws.Range("BP28").Value = "P"
$ws.Range("BP29").Value = "P"

# Mirror column BO's cell formatting (date format on row 1, centered text
# style on the rest) onto the new column BP. Row 12's player already left
# the team before this column range (no BO12 cell either), so it is
# deliberately excluded to avoid manufacturing a stray empty BP12 cell.
$ws.Range("BO1:BO11").Copy()
$ws.Range("BP1:BP11").PasteSpecial(-4122)
$ws.Range("BO13:BO29").Copy()
$ws.Range("BP13:BP29").PasteSpecial(-4122)

# Move the active selection, matching where the user last clicked after
# entering the new column of data.
$ws.Range("BS24").Select() | Out-Null
